# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) of the "Session Analysis Results" sheet.
#   "System, dnasr281@gmail.com"  ->  "dnasr281@gmail.com, System"
#
# Only cells whose current value is exactly the old ordering are touched,
# so already-correct / differently-valued cells (e.g. plain
# "dnasr281@gmail.com" or plain "System") are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

# Column G is the 7th column; data runs from row 2 through the last used row.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
